# Insert a new daily price record at row 154 for
# "Feria Lagunitas de Puerto Montt - Pepino ensalada", pushing the
# existing rows 154:193 down to 155:194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 154 downward (xlShiftDown semantics) to make room for the
# newly reported week's data.
$ws.Rows("154:154").Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A154").Value = 4
$ws.Range("B154").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C154").Value = "Los Lagos"
$ws.Range("D154").Value = 44543
$ws.Range("E154").Value = 10
$ws.Range("F154").Value = 100112043
$ws.Range("G154").Value = "Pepino ensalada"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 120
$ws.Range("K154").Value = 11000
$ws.Range("L154").Value = 11000
$ws.Range("M154").Value = 11000
$ws.Range("N154").Value = "`$/caja 60 unidades"
$ws.Range("O154").Value = "Región de Arica y Parinacota"
$ws.Range("P154").Value = 183
$ws.Range("Q154").Value = 60
$ws.Range("R154").Value = "Hortaliza"
